# Applies the diff: insert one new row (A 45895-2023) at row 11, bump the
# "Förändrad" (column C) date to 45202 on every existing data row, and
# append three brand-new rows (A 45529-2023, A 45871-2023, A 46436-2023)
# at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at position 11 ------------------------------
# This pushes the former rows 11-57 down to 12-58, copying formatting
# (date format on B/C, wrap text on R) from the row above, exactly like
# Excel's native "Insert Row" behaviour.
$ws.Rows.Item(11).Insert()

# --- Step 2: bump column C ("Förändrad") to 45202 on every data row -------
for ($r = 2; $r -le 58; $r++) {
    $ws.Range("C$r").Value2 = 45202
}

# --- Step 3: populate the newly inserted row 11 (A 45895-2023) ------------
$ws.Range("A11").Value2 = "A 45895-2023"
$ws.Range("B11").Value2 = 45195
$ws.Range("D11").Value2 = "UPPSALA LÄN"
$ws.Range("E11").Value2 = "ÄLVKARLEBY"
$ws.Range("G11").Value2 = 6.7
$ws.Range("H11").Value2 = 0
$ws.Range("I11").Value2 = 2
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = 0
$ws.Range("L11").Value2 = 0
$ws.Range("M11").Value2 = 0
$ws.Range("N11").Value2 = 0
$ws.Range("O11").Value2 = 0
$ws.Range("P11").Value2 = 0
$ws.Range("Q11").Value2 = 2
$ws.Range("R11").Value2 = "Sårläka`r`nTrådfräken"
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/artfynd/A 45895-2023.xlsx", "A 45895-2023")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/kartor/A 45895-2023.png", "A 45895-2023")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomål/A 45895-2023.docx", "A 45895-2023")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/klagomålsmail/A 45895-2023.docx", "A 45895-2023")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsyn/A 45895-2023.docx", "A 45895-2023")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALVKARLEBY/tillsynsmail/A 45895-2023.docx", "A 45895-2023")'
# wrapped 2-line Artnamn auto-grows the row height; pin it back to 15pt to
# match the sheet's standard row height.
$ws.Rows.Item(11).RowHeight = 15

# --- Step 4: the former last row (57) is now row 58 and should carry the
# same explicit 15pt custom height as every other interior row -----------
$ws.Rows.Item(58).RowHeight = 15

# --- Step 5: append brand-new row 59 (A 45529-2023) -----------------------
$ws.Range("A59").Value2 = "A 45529-2023"
$ws.Range("B59").Value2 = 45194
$ws.Range("C59").Value2 = 45202
$ws.Range("D59").Value2 = "UPPSALA LÄN"
$ws.Range("E59").Value2 = "ÄLVKARLEBY"
$ws.Range("F59").Value2 = "Bergvik skog väst AB"
$ws.Range("G59").Value2 = 2.6
$ws.Range("H59").Value2 = 0
$ws.Range("I59").Value2 = 0
$ws.Range("J59").Value2 = 0
$ws.Range("K59").Value2 = 0
$ws.Range("L59").Value2 = 0
$ws.Range("M59").Value2 = 0
$ws.Range("N59").Value2 = 0
$ws.Range("O59").Value2 = 0
$ws.Range("P59").Value2 = 0
$ws.Range("Q59").Value2 = 0
$ws.Range("B59:C59").NumberFormat = "YYYY-MM-DD"
$ws.Range("R59").WrapText = $true
$ws.Rows.Item(59).RowHeight = 15

# --- Step 6: append brand-new row 60 (A 45871-2023) -----------------------
$ws.Range("A60").Value2 = "A 45871-2023"
$ws.Range("B60").Value2 = 45195
$ws.Range("C60").Value2 = 45202
$ws.Range("D60").Value2 = "UPPSALA LÄN"
$ws.Range("E60").Value2 = "ÄLVKARLEBY"
$ws.Range("G60").Value2 = 3.7
$ws.Range("H60").Value2 = 0
$ws.Range("I60").Value2 = 0
$ws.Range("J60").Value2 = 0
$ws.Range("K60").Value2 = 0
$ws.Range("L60").Value2 = 0
$ws.Range("M60").Value2 = 0
$ws.Range("N60").Value2 = 0
$ws.Range("O60").Value2 = 0
$ws.Range("P60").Value2 = 0
$ws.Range("Q60").Value2 = 0
$ws.Range("B60:C60").NumberFormat = "YYYY-MM-DD"
$ws.Range("R60").WrapText = $true
$ws.Rows.Item(60).RowHeight = 15

# --- Step 7: append brand-new row 61 (A 46436-2023) -----------------------
# (left without an explicit custom row height, matching the diff)
$ws.Range("A61").Value2 = "A 46436-2023"
$ws.Range("B61").Value2 = 45197
$ws.Range("C61").Value2 = 45202
$ws.Range("D61").Value2 = "UPPSALA LÄN"
$ws.Range("E61").Value2 = "ÄLVKARLEBY"
$ws.Range("F61").Value2 = "Bergvik skog väst AB"
$ws.Range("G61").Value2 = 0.8
$ws.Range("H61").Value2 = 0
$ws.Range("I61").Value2 = 0
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 0
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = 0
$ws.Range("N61").Value2 = 0
$ws.Range("O61").Value2 = 0
$ws.Range("P61").Value2 = 0
$ws.Range("Q61").Value2 = 0
$ws.Range("B61:C61").NumberFormat = "YYYY-MM-DD"
$ws.Range("R61").WrapText = $true

Write-Output "Applied edits: inserted row 11, refreshed column C, appended rows 59-61."
